$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append row 12 (2021 data) below the existing 2011-2020 rows (rows 2-11).
# A12 needs the same "year label" style as the other column-A cells, so
# copy row 11's A cell (value + format) first, then overwrite the value.
$ws.Range("A11").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "2021年"

$ws.Range("B12").Value = 4815

# C12, L12, R12 and S12 are blank-but-text cells in this table (same as the
# corresponding cells in every other row). Writing "" via .Value clears a
# cell entirely in this engine, so force text first with a bare leading
# apostrophe, then re-copy the already-blank source cell on top to drop the
# incidental quote-prefix formatting that introduces.
$ws.Range("C12").Value = "'"
$ws.Range("C11").Copy($ws.Range("C12"))

$ws.Range("D12").Value = 11313
$ws.Range("E12").Value = 11609
$ws.Range("F12").Value = 24485
$ws.Range("G12").Value = 65699
$ws.Range("H12").Value = 31497
$ws.Range("I12").Value = 10376
$ws.Range("J12").Value = 5051
$ws.Range("K12").Value = 6085

$ws.Range("L12").Value = "'"
$ws.Range("L11").Copy($ws.Range("L12"))

$ws.Range("M12").Value = 265530
$ws.Range("N12").Value = 69755
$ws.Range("O12").Value = 22544
$ws.Range("P12").Value = 6366
$ws.Range("Q12").Value = 4779

$ws.Range("R12").Value = "'"
$ws.Range("R11").Copy($ws.Range("R12"))

$ws.Range("S12").Value = "'"
$ws.Range("S11").Copy($ws.Range("S12"))

$ws.Range("T12").Value = 1300
$ws.Range("U12").Value = 397524
